$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    3  = 500
    4  = 10000
    5  = 5000
    6  = 1000
    7  = 1500
    8  = 300
    10 = 500
    11 = 500
    13 = 500
    14 = 2000
    15 = 3000
    16 = 5000
    18 = 500
    19 = 1500
    24 = 15000
    25 = 500
    27 = 500
    28 = 5000
    29 = 300
    30 = 500
    31 = 2000
    33 = 8000
    34 = 5000
    35 = 100
    38 = 2000
    39 = 500
    40 = 1500
    44 = 500
    45 = 500
    46 = 500
    47 = 100
    50 = 500
    52 = 74800
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
